$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 19:04"

# Update country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) for the rows whose
# figures changed in today's data refresh.
$ws.Range("B4").Value = 1271059
$ws.Range("C4").Value = 7967
$ws.Range("D4").Value = 213562
$ws.Range("E4").Value = 981939
$ws.Range("F4").Value = 15827
$ws.Range("G4").Value = 759
$ws.Range("H4").Value = 75558
$ws.Range("B7").Value = 206715
$ws.Range("C7").Value = 5614
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 175756
$ws.Range("F7").Value = 1559
$ws.Range("G7").Value = 539
$ws.Range("H7").Value = 30615
$ws.Range("B10").Value = 168912
$ws.Range("C10").Value = 750
$ws.Range("D10").Value = 139900
$ws.Range("E10").Value = 21676
$ws.Range("F10").Value = 1823
$ws.Range("G10").Value = 61
$ws.Range("H10").Value = 7336
$ws.Range("B11").Value = 133721
$ws.Range("C11").Value = 1977
$ws.Range("D11").Value = 82984
$ws.Range("E11").Value = 47096
$ws.Range("F11").Value = 1260
$ws.Range("G11").Value = 57
$ws.Range("H11").Value = 3641
$ws.Range("B21").Value = 30126
$ws.Range("C21").Value = 66
$ws.Range("D21").Value = 25700
$ws.Range("E21").Value = 2616
$ws.Range("F21").Value = 121
$ws.Range("G21").Value = 5
$ws.Range("H21").Value = 1810
$ws.Range("B37").Value = 14499
$ws.Range("C37").Value = 392
$ws.Range("D37").Value = 6144
$ws.Range("E37").Value = 7469
$ws.Range("F37").Value = 234
$ws.Range("G37").Value = 22
$ws.Range("H37").Value = 886
$ws.Range("B43").Value = 10083
$ws.Range("C43").Value = 145
$ws.Range("D43").Value = 7711
$ws.Range("E43").Value = 1858
$ws.Range("F43").Value = 43
$ws.Range("G43").Value = 8
$ws.Range("H43").Value = 514
$ws.Range("B47").Value = 8002
$ws.Range("C47").Value = 28
$ws.Range("D47").Value = 4369
$ws.Range("E47").Value = 3364
$ws.Range("F47").Value = 52
$ws.Range("G47").Value = 7
$ws.Range("H47").Value = 269
$ws.Range("B48").Value = 7996
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 32
$ws.Range("E48").Value = 7748
$ws.Range("F48").Value = 27
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 216
$ws.Range("B49").Value = 7981
$ws.Range("C49").Value = 393
$ws.Range("D49").Value = 1887
$ws.Range("E49").Value = 5612
$ws.Range("F49").Value = 41
$ws.Range("G49").Value = 13
$ws.Range("H49").Value = 482
$ws.Range("B50").Value = 7808
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 3153
$ws.Range("E50").Value = 4502
$ws.Range("F50").Value = 36
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 153
$ws.Range("B51").Value = 7731
$ws.Range("C51").Value = 208
$ws.Range("D51").Value = 859
$ws.Range("E51").Value = 6654
$ws.Range("F51").Value = 88
$ws.Range("G51").Value = 8
$ws.Range("H51").Value = 218
$ws.Range("B56").Value = 5548
$ws.Range("C56").Value = 140
$ws.Range("D56").Value = 2179
$ws.Range("E56").Value = 3186
$ws.Range("F56").Value = 1
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 183
$ws.Range("B60").Value = 4530
$ws.Range("C60").Value = 108
$ws.Range("D60").Value = 1518
$ws.Range("E60").Value = 2982
$ws.Range("F60").Value = 31
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 30
$ws.Range("B62").Value = 3859
$ws.Range("C62").Value = 8
$ws.Range("D62").Value = 3505
$ws.Range("E62").Value = 254
$ws.Range("F62").Value = 21
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 100
$ws.Range("B72").Value = 2269
$ws.Range("C72").Value = 36
$ws.Range("D72").Value = 1656
$ws.Range("E72").Value = 603
$ws.Range("F72").Value = 8
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 10
$ws.Range("B95").Value = 928
$ws.Range("C95").Value = 55
$ws.Range("D95").Value = 106
$ws.Range("E95").Value = 778
$ws.Range("F95").Value = 2
$ws.Range("G95").Value = 5
$ws.Range("H95").Value = 44
$ws.Range("B96").Value = 909
$ws.Range("C96").Value = 9
$ws.Range("D96").Value = 464
$ws.Range("E96").Value = 427
$ws.Range("F96").Value = 3
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 18
$ws.Range("B97").Value = 895
$ws.Range("C97").Value = 24
$ws.Range("D97").Value = 637
$ws.Range("E97").Value = 246
$ws.Range("F97").Value = 13
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 12
$ws.Range("B98").Value = 889
$ws.Range("C98").Value = 6
$ws.Range("D98").Value = 296
$ws.Range("E98").Value = 578
$ws.Range("F98").Value = 15
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 15
$ws.Range("B146").Value = 153
$ws.Range("C146").Value = 7
$ws.Range("D146").Value = 103
$ws.Range("E146").Value = 46
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 4
$ws.Range("B147").Value = 152
$ws.Range("C147").Value = 0
$ws.Range("D147").Value = 104
$ws.Range("E147").Value = 35
$ws.Range("F147").Value = 4
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 13
$ws.Range("B205").Value = 11
$ws.Range("C205").Value = 0
$ws.Range("D205").Value = 8
$ws.Range("E205").Value = 3
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("B206").Value = 11
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 7
$ws.Range("E206").Value = 3
$ws.Range("F206").Value = 1
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1
